# Actualización automática del mapa (2025-07-23 07:17:33)
# Inserts a new incident row at row 52 (pushing existing rows 52-69 down to 53-70)
# and populates it with the new PEBCOM case data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 52..69 down to 53..70, creating a fresh blank row 52.
$ws.Rows.Item(52).Insert()

# Helper to write a value that must be stored as TEXT even though it
# looks like a number or a date (mirrors the source data, which keeps
# every one of these columns as plain text).
function Set-TextValue($cell, $value) {
    $cell.NumberFormat = "@"
    $cell.Value = $value
}

Set-TextValue $ws.Cells.Item(52, 1) "807044200"
Set-TextValue $ws.Cells.Item(52, 2) "5/29/2025"
$ws.Cells.Item(52, 3).Value = "11 de Septiembre de 1888 4662"
Set-TextValue $ws.Cells.Item(52, 4) "13"
Set-TextValue $ws.Cells.Item(52, 5) "807044200"
$ws.Cells.Item(52, 6).Value = "PEBCOM"
$ws.Cells.Item(52, 7).Value = "Pendiente"
$ws.Cells.Item(52, 8).Value = "CAMBIAR COLUMNA MUY INCLINADA POR POSTE PRFV 400, COLOCAR A 40 CMTS DEL CORDON"
$ws.Cells.Item(52, 9).Value = 1
$ws.Cells.Item(52, 10).Value = "Aplomo"
$ws.Cells.Item(52, 11).Value = "Sin equipos"
$ws.Cells.Item(52, 12).Value = "Terminal"
$ws.Cells.Item(52, 13).Value = -58.467458
$ws.Cells.Item(52, 14).Value = -34.537549
$ws.Cells.Item(52, 15).Value = "Saavedra"
$ws.Cells.Item(52, 16).Value = "Capital Norte"
